$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay stored as text, matching original inline-string formatting.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '90.615.76'
$ws.Range("E2").Value = '  -0.58%  '

# Row 3
$ws.Range("D3").Value = '3.141.88'
$ws.Range("E3").Value = '  +1.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").Value = '237.32'
$ws.Range("E5").Value = '  +8.14%  '

# Row 6
$ws.Range("D6").Value = '641.23'
$ws.Range("E6").Value = '  +2.89%  '

# Row 7
$ws.Range("D7").Value = '1.07'
$ws.Range("E7").Value = '  +11.41%  '

# Row 8
$ws.Range("E8").Value = '  -5.37%  '

# Row 9
$ws.Range("E9").Value = '  +0.12%  '

# Row 10
$ws.Range("D10").Value = '3.139.15'
$ws.Range("E10").Value = '  +1.16%  '

# Row 11
$ws.Range("D11").Value = '0.724'
$ws.Range("E11").Value = '  +0.38%  '

# Row 12
$ws.Range("E12").Value = '  +3.96%  '

# Row 13
$ws.Range("D13").Value = '36.60'
$ws.Range("E13").Value = '  +6.58%  '

# Row 14
$ws.Range("E14").Value = '  -4.73%  '

# Row 15
$ws.Range("D15").Value = '5.64'
$ws.Range("E15").Value = '  +4.52%  '

# Row 16
$ws.Range("D16").Value = '90.304.31'
$ws.Range("E16").Value = '  -0.74%  '

# Row 17
$ws.Range("D17").Value = '3.717.01'
$ws.Range("E17").Value = '  +0.83%  '

# Row 18
$ws.Range("D18").Value = '3.242.83'
$ws.Range("E18").Value = '  +3.20%  '

# Row 19
$ws.Range("D19").Value = '3.73'
$ws.Range("E19").Value = '  -0.39%  '

# Row 20
$ws.Range("B20").Value = 'PEPE'
$ws.Range("C20").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D20").Value = '0.0000219'
$ws.Range("E20").Value = '  -0.99%  '

# Row 21
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '14.48'
$ws.Range("E21").Value = '  +3.22%  '

# Row 22
$ws.Range("D22").Value = '450.74'
$ws.Range("E22").Value = '  +3.65%  '

# Row 23
$ws.Range("E23").Value = '  +9.98%  '

# Row 24
$ws.Range("E24").Value = '  +3.06%  '

# Row 25
$ws.Range("E25").Value = '  -2.49%  '

# Row 26
$ws.Range("D26").Value = '91.13'
$ws.Range("E26").Value = '  +4.29%  '

# Row 27
$ws.Range("D27").Value = '12.54'
$ws.Range("E27").Value = '  +3.02%  '

# Row 28
$ws.Range("E28").Value = '  +0.47%  '

# Row 30
$ws.Range("D30").Value = '9.91'
$ws.Range("E30").Value = '  +8.88%  '

# Row 31
$ws.Range("E31").Value = '  -4.34%  '

# Row 32
$ws.Range("D32").Value = '27.42'
$ws.Range("E32").Value = '  +15.91%  '

# Row 33
$ws.Range("D33").Value = '0.200'
$ws.Range("E33").Value = '  +33.02%  '

# Row 34
$ws.Range("D34").Value = '3.91'
$ws.Range("E34").Value = '  +4.21%  '

# Row 35
$ws.Range("D35").Value = '518.92'
$ws.Range("E35").Value = '  -0.97%  '

# Row 36
$ws.Range("D36").Value = '0.150'
$ws.Range("E36").Value = '  +5.22%  '

# Row 37
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  +5.13%  '

# Row 38
$ws.Range("D38").Value = '7.12'
$ws.Range("E38").Value = '  +0.44%  '

# Row 39
$ws.Range("E39").Value = '  +2.74%  '

# Row 40
$ws.Range("E40").Value = '  +8.27%  '

# Row 41
$ws.Range("D41").Value = '22.21'
$ws.Range("E41").Value = '  -0.29%  '

# Row 42
$ws.Range("D42").Value = '0.0863'
$ws.Range("E42").Value = '  +1.93%  '

# Row 44
$ws.Range("D44").Value = '0.744'
$ws.Range("E44").Value = '  -16.32%  '

# Row 45
$ws.Range("D45").Value = '3.38'
$ws.Range("E45").Value = '  +41.38%  '

# Row 46
$ws.Range("E46").Value = '  +2.06%  '

# Row 47
$ws.Range("D47").Value = '0.706'
$ws.Range("E47").Value = '  +14.21%  '

# Row 48
$ws.Range("D48").Value = '150.51'
$ws.Range("E48").Value = '  +2.36%  '

# Row 49
$ws.Range("D49").Value = '4.58'
$ws.Range("E49").Value = '  +9.77%  '

# Row 50
$ws.Range("D50").Value = '45.63'
$ws.Range("E50").Value = '  +3.62%  '

# Row 51
$ws.Range("E51").Value = '  +4.71%  '
